$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price ("D") and 1h volume-change ("E") cells for rows 2-51.
# New values are written with a leading apostrophe, exactly as a user
# would type a numeric-looking value into Excel to force it to stay
# literal text (e.g. '232.03', '0.620') instead of being auto-converted
# to a number and losing the trailing zero / thousands-dot formatting.
$updates = [ordered]@{
    'D2' = "'34.877.06"
    'E2' = "'  -0.71%  "
    'D3' = "'1.840.94"
    'E3' = "'  +1.52%  "
    'E4' = "'  +0.04%  "
    'D5' = "'232.03"
    'E5' = "'  -0.57%  "
    'D6' = "'0.620"
    'E6' = "'  +1.19%  "
    'E7' = "'  +0.05%  "
    'D8' = "'39.81"
    'E8' = "'  -2.02%  "
    'E9' = "'  +1.46%  "
    'D10' = "'0.0687"
    'E10' = "'  +0.41%  "
    'D11' = "'0.0984"
    'E11' = "'  -1.26%  "
    'D12' = "'2.107.15"
    'D13' = "'11.42"
    'E13' = "'  +3.37%  "
    'D14' = "'1.849.70"
    'E14' = "'  +1.85%  "
    'D15' = "'0.673"
    'E15' = "'  +1.38%  "
    'E16' = "'  -0.30%  "
    'D17' = "'34.931.10"
    'E17' = "'  -0.45%  "
    'D18' = "'69.91"
    'E18' = "'  +0.41%  "
    'E19' = "'  -0.37%  "
    'D20' = "'240.70"
    'E20' = "'  +0.79%  "
    'E21' = "'  +2.35%  "
    'D22' = "'4.68"
    'E22' = "'  -0.51%  "
    'E23' = "'  -0.06%  "
    'E24' = "'  +1.09%  "
    'D25' = "'171.62"
    'E25' = "'  -0.15%  "
    'D26' = "'7.79"
    'E26' = "'  -0.66%  "
    'D27' = "'17.45"
    'E27' = "'  -0.31%  "
    'E28' = "'  +2.31%  "
    'E29' = "'  -5.78%  "
    'E30' = "'  +0.07%  "
    'E31' = "'  +0.08%  "
    'E32' = "'  -5.55%  "
    'E33' = "'  -1.65%  "
    'E34' = "'  +8.31%  "
    'E35' = "'  +8.56%  "
    'E36' = "'  +11.79%  "
    'D37' = "'0.696"
    'E37' = "'  +2.65%  "
    'E38' = "'  +6.91%  "
    'D39' = "'91.09"
    'E39' = "'  -1.46%  "
    'D40' = "'1.345.25"
    'E40' = "'  +2.44%  "
    'E41' = "'  +0.35%  "
    'D42' = "'14.88"
    'E42' = "'  +2.17%  "
    'E43' = "'  -0.11%  "
    'E44' = "'  -2.71%  "
    'E45' = "'  -0.11%  "
    'E46' = "'  -0.05%  "
    'E47' = "'  +2.12%  "
    'D48' = "'2.020.87"
    'E48' = "'  +1.49%  "
    'D49' = "'3.40"
    'E49' = "'  +19.65%  "
    'E50' = "'  +0.03%  "
    'D51' = "'0.0664"
    'E51' = "'  +1.92%  "
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
